$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 144.35
$ws.Range("I15").Value = 144.35
$ws.Range("K15").Value = 433.05
$ws.Range("M15").Value = -264.05
$ws.Range("H100").Value = 1878.2858
$ws.Range("I100").Value = 1204.1538
$ws.Range("J100").Value = 2973.75
$ws.Range("K100").Value = 1204.1538
$ws.Range("L100").Value = 2973.75
$ws.Range("M100").Value = -663.1538
$ws.Range("N100").Value = -4055.75
$ws.Range("H112").Value = 1337850
$ws.Range("I112").Value = 389.66666
$ws.Range("J112").Value = 1624448.8
$ws.Range("K112").Value = 1168.99998
$ws.Range("L112").Value = 4873346.4
$ws.Range("M112").Value = -60.99998000000005
$ws.Range("N112").Value = -4875562.4
$ws.Range("H127").Value = 961.3684
$ws.Range("I127").Value = 680.5833
$ws.Range("J127").Value = 1442.7142
$ws.Range("K127").Value = 2041.7499
$ws.Range("L127").Value = 4328.142599999999
$ws.Range("M127").Value = 2918.2501
$ws.Range("N127").Value = -14248.1426
$ws.Range("H129").Value = 6053.5713
$ws.Range("I129").Value = 386.22223
$ws.Range("J129").Value = 7599.212
$ws.Range("K129").Value = 1158.66669
$ws.Range("L129").Value = 22797.636
$ws.Range("M129").Value = 3841.33331
$ws.Range("N129").Value = -32797.636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 28183.334
$ws.Range("J92").Value = 28183.334
$ws.Range("L92").Value = 28183.334
$ws.Range("N92").Value = -33175.334
$ws.Range("H122").Value = 3127.5483
$ws.Range("I122").Value = 2745.4546
$ws.Range("J122").Value = 4061.5557
$ws.Range("K122").Value = 8236.363799999999
$ws.Range("L122").Value = 12184.6671
$ws.Range("M122").Value = -5786.363799999999
$ws.Range("N122").Value = -17084.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = ""
$ws.Range("H38").Value = 9999
$ws.Range("J38").Value = 9999
$ws.Range("L38").Value = 9999
$ws.Range("N38").Value = -10831
$ws.Range("H92").Value = 41200.5
$ws.Range("J92").Value = 41200.5
$ws.Range("L92").Value = 41200.5
$ws.Range("N92").Value = -46192.5
$ws.Range("H99").Value = 2016
$ws.Range("I99").Value = 1808.3846
$ws.Range("J99").Value = 2315.889
$ws.Range("K99").Value = 1808.3846
$ws.Range("L99").Value = 2315.889
$ws.Range("M99").Value = -310.3846000000001
$ws.Range("N99").Value = -5311.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2374.075
$ws.Range("I31").Value = 1593.3334
$ws.Range("J31").Value = 3012.8635
$ws.Range("K31").Value = 1593.3334
$ws.Range("L31").Value = 3012.8635
$ws.Range("M31").Value = -1298.3334
$ws.Range("N31").Value = -3602.8635
$ws.Range("H34").Value = 2374.075
$ws.Range("I34").Value = 1593.3334
$ws.Range("J34").Value = 3012.8635
$ws.Range("K34").Value = 1593.3334
$ws.Range("L34").Value = 3012.8635
$ws.Range("M34").Value = -1391.3334
$ws.Range("N34").Value = -3416.8635
$ws.Range("H35").Value = 1925
$ws.Range("I35").Value = 1925
$ws.Range("K35").Value = 1925
$ws.Range("M35").Value = -1631
$ws.Range("H58").Value = 2299.5366
$ws.Range("I58").Value = 825.125
$ws.Range("J58").Value = 4381.0586
$ws.Range("K58").Value = 825.125
$ws.Range("L58").Value = 4381.0586
$ws.Range("M58").Value = -622.125
$ws.Range("N58").Value = -4787.0586
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H125").Value = 49698
$ws.Range("J125").Value = 49698
$ws.Range("L125").Value = 49698
$ws.Range("N125").Value = -54618
$ws.Range("H136").Value = 2299.5366
$ws.Range("I136").Value = 825.125
$ws.Range("J136").Value = 4381.0586
$ws.Range("K136").Value = 2475.375
$ws.Range("L136").Value = 13143.1758
$ws.Range("M136").Value = 74.625
$ws.Range("N136").Value = -18243.1758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 797.3889
$ws.Range("I68").Value = 435
$ws.Range("J68").Value = 936.7692
$ws.Range("K68").Value = 1305
$ws.Range("L68").Value = 2810.3076
$ws.Range("M68").Value = -494
$ws.Range("N68").Value = -4432.3076
$ws.Range("H71").Value = 797.3889
$ws.Range("I71").Value = 435
$ws.Range("J71").Value = 936.7692
$ws.Range("K71").Value = 3915
$ws.Range("L71").Value = 8430.9228
$ws.Range("M71").Value = 141
$ws.Range("N71").Value = -16542.9228
$ws.Range("H131").Value = 1755243.9
$ws.Range("I131").Value = 7407834
$ws.Range("J131").Value = 991.8276
$ws.Range("K131").Value = 22223502
$ws.Range("L131").Value = 2975.4828
$ws.Range("M131").Value = -22218462
$ws.Range("N131").Value = -13055.4828
$ws.Range("H132").Value = 1237.375
$ws.Range("I132").Value = 518.75
$ws.Range("J132").Value = 1596.6875
$ws.Range("K132").Value = 4668.75
$ws.Range("L132").Value = 14370.1875
$ws.Range("M132").Value = -2138.75
$ws.Range("N132").Value = -19430.1875
$ws.Range("H134").Value = 791.5
$ws.Range("I134").Value = 473.6875
$ws.Range("K134").Value = 1421.0625
$ws.Range("M134").Value = 3648.9375
$ws.Range("H137").Value = 2950.2424
$ws.Range("I137").Value = 895.4375
$ws.Range("J137").Value = 4884.1763
$ws.Range("K137").Value = 2686.3125
$ws.Range("L137").Value = 14652.5289
$ws.Range("M137").Value = 2413.6875
$ws.Range("N137").Value = -24852.5289
$ws.Range("H139").Value = 3922.6843
$ws.Range("I139").Value = 1263.8334
$ws.Range("J139").Value = 8480.714
$ws.Range("K139").Value = 3791.5002
$ws.Range("L139").Value = 25442.142
$ws.Range("M139").Value = 1348.4998
$ws.Range("N139").Value = -35722.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 3188
$ws.Range("I32").Value = 1025.6
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 1025.6
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -708.5999999999999
$ws.Range("N32").Value = -14634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2273.65
$ws.Range("I126").Value = 2769.6428
$ws.Range("J126").Value = 1116.3334
$ws.Range("K126").Value = 8308.928400000001
$ws.Range("L126").Value = 3349.0002
$ws.Range("M126").Value = -5838.928400000001
$ws.Range("N126").Value = -8289.0002

Write-Host "Applied 162 cell updates across 8 sheets."
